# Doc fix by konghao
# Fill in the "备注" (remark) column for rows 12-27 of the call-log sheet
# with the updated [创达] test-precondition notes (3-record variants),
# replacing the earlier 1-record-variant notes that were never filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(12, 4).Value = "[创达]前提：手机中有已存的通话记录，其中包含已接电话3条、未接电话3条、已拨电话3条`n1.选中全部通话条目`nAssert1：页面上有9条通话记录"
$ws.Cells.Item(13, 4).Value = "[创达]前提：手机中有已存的通话记录，其中包含已接电话3条、未接电话3条、已拨电话3条`n1.选中未接电话条目`nAssert1：页面上有3条未接电话通话记录"
$ws.Cells.Item(14, 4).Value = "[创达]前提：手机中有已存的通话记录，其中包含已接电话3条、未接电话3条、已拨电话3条`n1.选中已拨电话电话条目`nAssert1：页面上有3条已拨电话通话记录"
$ws.Cells.Item(15, 4).Value = "[创达]前提：手机中有已存的通话记录，其中包含已接电话3条、未接电话3条、已拨电话3条`n1.选中已接电话电话条目`nAssert1：页面上有3条已接电话通话记录"
$ws.Cells.Item(16, 4).Value = "[创达]前提：手机中有已存的通话记录，3条今天（手机时间上是今天）的通话记录"
$ws.Cells.Item(17, 4).Value = "[创达]前提：手机中有已存的通话记录，3条昨天（手机时间上今天，相较于昨天是今天）的通话记录"
$ws.Cells.Item(18, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(19, 4).Value = "[创达]前提：手机中有已存的通话记录，3条今天（手机时间上是今天）的通话记录"
$ws.Cells.Item(20, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(21, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(22, 4).Value = "[创达]前提：手机中有已存的通话记录，3条今天（手机时间上是今天）的通话记录"
$ws.Cells.Item(23, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(24, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(25, 4).Value = "[创达]前提：手机中有已存的通话记录，3条今天（手机时间上是今天）的通话记录"
$ws.Cells.Item(26, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"
$ws.Cells.Item(27, 4).Value = "[创达]前提：手机中有已存的通话记录，早于昨天的3条记录（手机时间上是今天，只要比当前时间早二天以上即可）的通话记录"

# Row 12's note grew to three lines; match the author's explicit row-height bump.
$ws.Rows.Item(12).RowHeight = 71.25

# Leave the view where the author left it: scrolled down to the edited rows,
# with the first edited cell selected.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()

